# Auto-generated Excel COM-interop script
# Updates market-price derived cells across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, WVR)
# to match the latest scheduled data pull, per commit 'chore: update Sheets via scheduled runner'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1828.5
$ws.Range("J17").Value = 1911.8572
$ws.Range("L17").Value = 5735.571599999999
$ws.Range("N17").Value = -6071.571599999999
$ws.Range("H43").Value = 1774.0952
$ws.Range("I43").Value = 1676.1818
$ws.Range("K43").Value = 1676.1818
$ws.Range("M43").Value = -1607.1818
$ws.Range("H55").Value = 374.125
$ws.Range("I55").Value = 428.1111
$ws.Range("K55").Value = 428.1111
$ws.Range("M55").Value = -214.1111
$ws.Range("H62").Value = 1732.875
$ws.Range("I62").Value = 1380.5714
$ws.Range("J62").Value = 4199
$ws.Range("K62").Value = 1380.5714
$ws.Range("L62").Value = 4199
$ws.Range("M62").Value = -756.5714
$ws.Range("N62").Value = -5447
$ws.Range("H64").Value = 90913460
$ws.Range("I64").Value = 4514.8335
$ws.Range("J64").Value = 200004180
$ws.Range("K64").Value = 4514.8335
$ws.Range("L64").Value = 200004180
$ws.Range("M64").Value = -4266.8335
$ws.Range("N64").Value = -200004676
$ws.Range("H65").Value = 1732.875
$ws.Range("I65").Value = 1380.5714
$ws.Range("J65").Value = 4199
$ws.Range("K65").Value = 6902.857
$ws.Range("L65").Value = 20995
$ws.Range("M65").Value = -3782.857
$ws.Range("N65").Value = -27235
$ws.Range("H67").Value = 90913460
$ws.Range("I67").Value = 4514.8335
$ws.Range("J67").Value = 200004180
$ws.Range("K67").Value = 4514.8335
$ws.Range("L67").Value = 200004180
$ws.Range("M67").Value = -3656.8335
$ws.Range("N67").Value = -200005896
$ws.Range("I100").Value = 1997
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1997
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1456
$ws.Range("N100").Value = ""
$ws.Range("H101").Value = 346.2
$ws.Range("J101").Value = 419.5
$ws.Range("L101").Value = 1258.5
$ws.Range("N101").Value = -4502.5
$ws.Range("H137").Value = 2383889.2
$ws.Range("I137").Value = 2942712.5
$ws.Range("K137").Value = 8828137.5
$ws.Range("M137").Value = -8825587.5
$ws.Range("H138").Value = 2954.0256
$ws.Range("J138").Value = 3287.3726
$ws.Range("L138").Value = 9862.1178
$ws.Range("N138").Value = -20142.1178

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1926004
$ws.Range("I32").Value = 883626.3
$ws.Range("K32").Value = 883626.3
$ws.Range("M32").Value = -883339.3
$ws.Range("H61").Value = 2449.6
$ws.Range("I61").Value = 1945.8572
$ws.Range("K61").Value = 1945.8572
$ws.Range("M61").Value = -1733.8572
$ws.Range("H63").Value = 1926
$ws.Range("J63").Value = 1999
$ws.Range("L63").Value = 1999
$ws.Range("N63").Value = -3371
$ws.Range("H66").Value = 1926
$ws.Range("J66").Value = 1999
$ws.Range("L66").Value = 9995
$ws.Range("N66").Value = -16859
$ws.Range("H74").Value = 2673.76
$ws.Range("I74").Value = 2461.7058
$ws.Range("J74").Value = 3124.375
$ws.Range("K74").Value = 2461.7058
$ws.Range("L74").Value = 3124.375
$ws.Range("M74").Value = -1587.7058
$ws.Range("N74").Value = -4872.375
$ws.Range("H77").Value = 2673.76
$ws.Range("I77").Value = 2461.7058
$ws.Range("J77").Value = 3124.375
$ws.Range("K77").Value = 12308.529
$ws.Range("L77").Value = 15621.875
$ws.Range("M77").Value = -7940.529
$ws.Range("N77").Value = -24357.875
$ws.Range("H97").Value = 1020
$ws.Range("I97").Value = 1020
$ws.Range("K97").Value = 1020
$ws.Range("M97").Value = -524
$ws.Range("H132").Value = 2840.4546
$ws.Range("I132").Value = 2499.3076
$ws.Range("K132").Value = 7497.9228
$ws.Range("M132").Value = -4967.9228
$ws.Range("H136").Value = 2449.6
$ws.Range("I136").Value = 1945.8572
$ws.Range("K136").Value = 5837.571599999999
$ws.Range("M136").Value = -3287.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17862262
$ws.Range("J20").Value = 3185.2
$ws.Range("L20").Value = 3185.2
$ws.Range("N20").Value = -3679.2
$ws.Range("H134").Value = 2935.8635
$ws.Range("I134").Value = 2238.3333
$ws.Range("K134").Value = 6714.999899999999
$ws.Range("M134").Value = -4179.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2362087
$ws.Range("I31").Value = 2882.8823
$ws.Range("K31").Value = 2882.8823
$ws.Range("M31").Value = -2587.8823
$ws.Range("H34").Value = 2362087
$ws.Range("I34").Value = 2882.8823
$ws.Range("K34").Value = 2882.8823
$ws.Range("M34").Value = -2680.8823
$ws.Range("H99").Value = 2997.25
$ws.Range("I99").Value = 1995
$ws.Range("K99").Value = 1995
$ws.Range("M99").Value = -497
$ws.Range("H126").Value = 2997.25
$ws.Range("I126").Value = 1995
$ws.Range("K126").Value = 5985
$ws.Range("M126").Value = -3515
$ws.Range("H132").Value = 16672522
$ws.Range("I132").Value = 3926.25
$ws.Range("J132").Value = 41675416
$ws.Range("K132").Value = 11778.75
$ws.Range("L132").Value = 125026248
$ws.Range("M132").Value = -9248.75
$ws.Range("N132").Value = -125031308

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 960293.3
$ws.Range("I4").Value = 262957.78
$ws.Range("K4").Value = 788873.3400000001
$ws.Range("M4").Value = -788761.3400000001
$ws.Range("H5").Value = 706.1429000000001
$ws.Range("I5").Value = 828
$ws.Range("K5").Value = 2484
$ws.Range("M5").Value = -2372
$ws.Range("H7").Value = 134.71428
$ws.Range("I7").Value = 134.71428
$ws.Range("K7").Value = 404.14284
$ws.Range("M7").Value = -292.14284
$ws.Range("H33").Value = 322.16666
$ws.Range("I33").Value = 405.66666
$ws.Range("J33").Value = 238.66667
$ws.Range("K33").Value = 2433.99996
$ws.Range("L33").Value = 1432.00002
$ws.Range("M33").Value = -2150.99996
$ws.Range("N33").Value = -1998.00002
$ws.Range("H68").Value = 12509038
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 14295943
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 42887829
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -42889451
$ws.Range("H71").Value = 12509038
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 14295943
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 128663487
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -128671599
$ws.Range("H135").Value = 706.1429000000001
$ws.Range("I135").Value = 828
$ws.Range("K135").Value = 7452
$ws.Range("M135").Value = -4917

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 45457980
$ws.Range("I80").Value = 66669676
$ws.Range("J80").Value = 4341.4287
$ws.Range("K80").Value = 66669676
$ws.Range("L80").Value = 4341.4287
$ws.Range("M80").Value = -66668678
$ws.Range("N80").Value = -6337.4287
$ws.Range("H83").Value = 45457980
$ws.Range("I83").Value = 66669676
$ws.Range("J83").Value = 4341.4287
$ws.Range("K83").Value = 333348380
$ws.Range("L83").Value = 21707.1435
$ws.Range("M83").Value = -333343388
$ws.Range("N83").Value = -31691.1435
$ws.Range("H132").Value = 2340.3784
$ws.Range("I132").Value = 2181.037
$ws.Range("K132").Value = 6543.110999999999
$ws.Range("M132").Value = -4013.110999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3134.5
$ws.Range("I81").Value = 2775.5
$ws.Range("K81").Value = 5551
$ws.Range("M81").Value = -4490
$ws.Range("H84").Value = 3134.5
$ws.Range("I84").Value = 2775.5
$ws.Range("K84").Value = 27755
$ws.Range("M84").Value = -22451
$ws.Range("H100").Value = 66668172
$ws.Range("I100").Value = 1760.5
$ws.Range("K100").Value = 3521
$ws.Range("M100").Value = -2980
